$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.230.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "'1.587.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'211.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'19.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'1.808.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'1.582.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'64.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "'26.217.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'0.0₃0725"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'215.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'4.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "'2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'8.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'143.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'6.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'15.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "'3.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'1.362.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'0.582"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'0.819"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'5.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("E43").Value = "  -17.27%  "
$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'1.721.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "'60.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "'86.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "'1.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "'0.0982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -1.09%  "
